$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 4 for columns D, J, K, L, M, O, P
$ws.Range("D2").Value = 44273
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 233

$ws.Range("D4").Value = 44291
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("O4").Value = "Limache"
$ws.Range("P4").Value = 183
